# Design add: jr and jal instructions
# Insert two new rows at row 39 (pushing the existing beq/bne/blez/bgtz rows
# down from 39-42 to 41-44) and populate them with the jal / jr truth-table
# entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 39 (shifts old rows 39..42 -> 41..44)
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(39).Insert()

# ---- Row 39: jal ----
$ws.Range("A39").Value2 = "jal"
$ws.Range("B39").Value2 = "000011"
$ws.Range("D39").Value2 = "jal_br"
$ws.Range("E39").Value2 = "x"
$ws.Range("F39").Value2 = "x"
$ws.Range("G39").Value2 = "x"
$ws.Range("H39").Value2 = "off"
$ws.Range("I39").Value2 = "x"
$ws.Range("J39").Value2 = "x"
$ws.Range("K39").Value2 = "x"
$ws.Range("L39").Value2 = "on"

# ---- Row 40: jr ----
$ws.Range("A40").Value2 = "jr"
$ws.Range("B40").Value2 = "000000"
$ws.Range("C40").Value2 = "001000"
$ws.Range("D40").Value2 = "jr_br"
$ws.Range("E40").Value2 = "x"
$ws.Range("F40").Value2 = "x"
$ws.Range("G40").Value2 = "x"
$ws.Range("H40").Value2 = "off"
$ws.Range("I40").Value2 = "x"
$ws.Range("J40").Value2 = "x"
$ws.Range("K40").Value2 = "x"
$ws.Range("L40").Value2 = "off"

# The "opcode" column (B) for these two new rows uses a text-formatted
# style without the extra left-alignment override that the rest of the
# table's B/C columns use (numFmtId 49 "@" + vertical-center only).
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").HorizontalAlignment = 1
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").HorizontalAlignment = 1

# Match the author's final cursor position / scroll so the saved view
# looks the same (active cell L39, scrolled down toward the bottom of
# the table).
$ws.Range("L39").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
